# FlightData.xlsx update:
#  - The "Departure Date" values on the Input sheet (25/01/2022) were no
#    longer valid test dates, so they were bumped to 25/02/2022.
#  - The Output sheet's most recent run (row 2) was refreshed with the
#    results of the latest test execution (new timestamp + fares).
#  - The last active cell/selection on the Input sheet moved to H14.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets("Input")
$wsOutput = $wb.Worksheets("Output")

# --- Input sheet: "Departure Date" column (F), rows 2, 4, 5 and 7
#     25012022 -> 25022022 (rows 3 and 6 were already 25022022)
$wsInput.Range("F2").Value = 25022022
$wsInput.Range("F4").Value = 25022022
$wsInput.Range("F5").Value = 25022022
$wsInput.Range("F7").Value = 25022022

# --- Output sheet: most recent test run entry (row 2) refreshed
$wsOutput.Range("A2").Value = "27/01/2022 10:27:17 am"
$wsOutput.Range("D2").Value = "₹1,21,080"
$wsOutput.Range("E2").Value = "₹8,344"
$wsOutput.Range("G2").Value = "₹1,29,434"

# --- Restore the last selection on the Input sheet (tab stays selected)
$wsInput.Activate()
$wsInput.Range("H14").Select()
